$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1139.4736
$ws.Range("J112").Value = 1138.8889
$ws.Range("L112").Value = 3416.6667
$ws.Range("N112").Value = -5632.6667

$ws.Range("H138").Value = 2878.7778
$ws.Range("I138").Value = 1416.7693
$ws.Range("K138").Value = 4250.3079
$ws.Range("M138").Value = 889.6921000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4691.815
$ws.Range("I32").Value = 4099.479
$ws.Range("K32").Value = 4099.479
$ws.Range("M32").Value = -3812.479

$ws.Range("H74").Value = 3140.2114
$ws.Range("I74").Value = 3413.8445
$ws.Range("J74").Value = 1381.1428
$ws.Range("K74").Value = 3413.8445
$ws.Range("L74").Value = 1381.1428
$ws.Range("M74").Value = -2539.8445
$ws.Range("N74").Value = -3129.1428

$ws.Range("H77").Value = 3140.2114
$ws.Range("I77").Value = 3413.8445
$ws.Range("J77").Value = 1381.1428
$ws.Range("K77").Value = 17069.2225
$ws.Range("L77").Value = 6905.714
$ws.Range("M77").Value = -12701.2225
$ws.Range("N77").Value = -15641.714

$ws.Range("H102").Value = 1244.375
$ws.Range("I102").Value = 1084.7727
$ws.Range("K102").Value = 1084.7727
$ws.Range("M102").Value = 537.2273

$ws.Range("H110").Value = 23278.125
$ws.Range("I110").Value = 29678.055
$ws.Range("K110").Value = 29678.055
$ws.Range("M110").Value = -27633.055

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 20000
$ws.Range("J31").Value = 20000
$ws.Range("L31").Value = 20000
$ws.Range("N31").Value = -20504

$ws.Range("H105").Value = 3829.5
$ws.Range("I105").Value = 3974.375
$ws.Range("K105").Value = 3974.375
$ws.Range("M105").Value = -2227.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3466.75
$ws.Range("I31").Value = 3183.3333
$ws.Range("J31").Value = 3601
$ws.Range("K31").Value = 3183.3333
$ws.Range("L31").Value = 3601
$ws.Range("M31").Value = -2888.3333
$ws.Range("N31").Value = -4191

$ws.Range("H34").Value = 3466.75
$ws.Range("I34").Value = 3183.3333
$ws.Range("J34").Value = 3601
$ws.Range("K34").Value = 3183.3333
$ws.Range("L34").Value = 3601
$ws.Range("M34").Value = -2981.3333
$ws.Range("N34").Value = -4005

$ws.Range("H45").Value = 3067
$ws.Range("I45").Value = 3067
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 3067
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2474
$ws.Range("N45").ClearContents()

$ws.Range("H51").Value = 32112.8
$ws.Range("J51").Value = 32112.8
$ws.Range("L51").Value = 32112.8
$ws.Range("N51").Value = -33584.8

$ws.Range("H58").Value = 1115.6757
$ws.Range("I58").Value = 723.2241
$ws.Range("J58").Value = 2538.3125
$ws.Range("K58").Value = 723.2241
$ws.Range("L58").Value = 2538.3125
$ws.Range("M58").Value = -520.2241
$ws.Range("N58").Value = -2944.3125

$ws.Range("H61").Value = 32112.8
$ws.Range("J61").Value = 32112.8
$ws.Range("L61").Value = 32112.8
$ws.Range("N61").Value = -32808.8

$ws.Range("H136").Value = 1115.6757
$ws.Range("I136").Value = 723.2241
$ws.Range("J136").Value = 2538.3125
$ws.Range("K136").Value = 2169.6723
$ws.Range("L136").Value = 7614.9375
$ws.Range("M136").Value = 380.3276999999998
$ws.Range("N136").Value = -12714.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 732101.0600000001
$ws.Range("I5").Value = 624.6667
$ws.Range("J5").Value = 1951228.4
$ws.Range("K5").Value = 1874.0001
$ws.Range("L5").Value = 5853685.199999999
$ws.Range("M5").Value = -1762.0001
$ws.Range("N5").Value = -5853909.199999999

$ws.Range("H131").Value = 6487
$ws.Range("I131").Value = 894.2857
$ws.Range("J131").Value = 9498.462
$ws.Range("K131").Value = 2682.8571
$ws.Range("L131").Value = 28495.386
$ws.Range("M131").Value = 2357.1429
$ws.Range("N131").Value = -38575.386

$ws.Range("H132").Value = 1761.25
$ws.Range("I132").Value = 1070
$ws.Range("J132").Value = 2326.818
$ws.Range("K132").Value = 9630
$ws.Range("L132").Value = 20941.362
$ws.Range("M132").Value = -7100
$ws.Range("N132").Value = -26001.362

$ws.Range("H135").Value = 732101.0600000001
$ws.Range("I135").Value = 624.6667
$ws.Range("J135").Value = 1951228.4
$ws.Range("K135").Value = 5622.0003
$ws.Range("L135").Value = 17561055.6
$ws.Range("M135").Value = -3087.0003
$ws.Range("N135").Value = -17566125.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1992.3889
$ws.Range("I80").Value = 1988.7778
$ws.Range("J80").Value = 1996
$ws.Range("K80").Value = 1988.7778
$ws.Range("L80").Value = 1996
$ws.Range("M80").Value = -990.7778000000001
$ws.Range("N80").Value = -3992

$ws.Range("H83").Value = 1992.3889
$ws.Range("I83").Value = 1988.7778
$ws.Range("J83").Value = 1996
$ws.Range("K83").Value = 9943.889000000001
$ws.Range("L83").Value = 9980
$ws.Range("M83").Value = -4951.889000000001
$ws.Range("N83").Value = -19964

$ws.Range("H132").Value = 2351.7036
$ws.Range("I132").Value = 2074.15
$ws.Range("J132").Value = 3144.7144
$ws.Range("K132").Value = 6222.450000000001
$ws.Range("L132").Value = 9434.143199999999
$ws.Range("M132").Value = -3692.450000000001
$ws.Range("N132").Value = -14494.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 4053.75
$ws.Range("I35").Value = 1238.3334
$ws.Range("K35").Value = 1238.3334
$ws.Range("M35").Value = -902.3334

$ws.Range("H40").Value = 77270.57000000001
$ws.Range("I40").Value = 104198.8
$ws.Range("J40").Value = 9950
$ws.Range("K40").Value = 104198.8
$ws.Range("L40").Value = 9950
$ws.Range("M40").Value = -104062.8
$ws.Range("N40").Value = -10222

$ws.Range("H93").Value = 4991.9165
$ws.Range("I93").Value = 6128.3335
$ws.Range("K93").Value = 6128.3335
$ws.Range("M93").Value = -4880.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 57526.277
$ws.Range("I122").Value = 92698.91
$ws.Range("J122").Value = 2255
$ws.Range("K122").Value = 278096.73
$ws.Range("L122").Value = 6765
$ws.Range("M122").Value = -275646.73
$ws.Range("N122").Value = -11665

$ws.Range("H132").Value = 1394.473
$ws.Range("I132").Value = 1162.2174
$ws.Range("J132").Value = 4599.6
$ws.Range("K132").Value = 3486.6522
$ws.Range("L132").Value = 13798.8
$ws.Range("M132").Value = -956.6522
$ws.Range("N132").Value = -18858.8
